# Generate Report for Handoff
# Refresh the localization-status report: flip the per-language status back
# to "Ready for handoff" and bump the handoff timestamps, then re-narrow the
# now-shorter "Status" columns on each sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Refreshed handoff-generation timestamps
$overview.Range("G2").Value = "2016-08-15 12:54:26"
$zhcn.Range("H2").Value     = "2016-08-15 12:54:22"
$dede.Range("H2").Value     = "2016-08-15 12:54:26"

# --- Narrower "Status" columns now that the text is shorter
$overview.Columns("E").ColumnWidth = 16.333333333333332
$overview.Columns("F").ColumnWidth = 16.333333333333332
$zhcn.Columns("C").ColumnWidth     = 16.333333333333332
$dede.Columns("C").ColumnWidth     = 16.333333333333332
